# Update "想去人数" (want-to-go count) values in the "F" column on the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets to match the
# refreshed data pulled from the source at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 318
$ws1.Range("F4").Value = 243
$ws1.Range("F5").Value = 2857
$ws1.Range("F6").Value = 1982
$ws1.Range("F7").Value = 381
$ws1.Range("F9").Value = 1047
$ws1.Range("F11").Value = 310
$ws1.Range("F12").Value = 45

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 318
$ws4.Range("F4").Value = 243
$ws4.Range("F5").Value = 2857
$ws4.Range("F6").Value = 1982
$ws4.Range("F7").Value = 381
$ws4.Range("F10").Value = 1047
$ws4.Range("F12").Value = 310
$ws4.Range("F13").Value = 45
